$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173013210296631
$ws.Range("B1").Value = 2.640596389770508
$ws.Range("C1").Value = 5.809912204742432
$ws.Range("D1").Value = 2.098309516906738
$ws.Range("E1").Value = 1.208070635795593
